# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet named "2022-Q1" right before the "总计" (total)
#    sheet and fill it with the per-fund holdings detail for 2022-Q1.
# 2. Insert a new top row into the "总计" summary sheet for the 2022-Q1
#    quarter, shifting the existing rows down, and renumber the index
#    column accordingly.

$wb = $excel.ActiveWorkbook

# --- Step 1: create & populate the "2022-Q1" detail sheet -----------------

$totalSheetBeforeInsert = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($totalSheetBeforeInsert, $null)
$newSheet.Name = "2022-Q1"

# Worksheet references are position-anchored, and inserting a sheet shifts
# every sheet at/after the insertion point -- re-resolve "总计" by name now
# that it has moved one slot to the right, otherwise $totalSheetBeforeInsert
# (and anything derived from it) would silently alias the new sheet instead.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy header row + data rows from an existing quarter sheet so the new
# sheet picks up identical styles (bold/centered/bordered header, bordered
# index column, etc.) without hand-rolling style indices. Column A of the
# header row is untouched in the source sheets (no cell there at all), so
# it is deliberately excluded to avoid materialising a spurious empty A1.
$template.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$template.Range("A2:H10").Copy($newSheet.Range("A2:H10"))

# Header text
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund code/name/scale/position/share/value are stored as text (inlineStr)
# in the source workbook, while the index (A) and rank (H) columns are
# numeric -- force the text columns' format so the values aren't
# reinterpreted as numbers.
$newSheet.Range("B2:G10").NumberFormat = "@"

$rows = New-Object 'object[,]' 9,8
$rows[0,0] = 0;  $rows[0,1] = "070099"; $rows[0,2] = "嘉实优质企业混合";          $rows[0,3] = "22.02"; $rows[0,4] = "91.84"; $rows[0,5] = "3.63"; $rows[0,6] = "0.7993"; $rows[0,7] = 9
$rows[1,0] = 1;  $rows[1,1] = "014207"; $rows[1,2] = "华安产业精选混合A";         $rows[1,3] = "27.31"; $rows[1,4] = "62.03"; $rows[1,5] = "2.03"; $rows[1,6] = "0.5544"; $rows[1,7] = 7
$rows[2,0] = 2;  $rows[2,1] = "010147"; $rows[2,2] = "博道嘉兴一年持有期混合";     $rows[2,3] = "14.00"; $rows[2,4] = "92.27"; $rows[2,5] = "3.58"; $rows[2,6] = "0.5012"; $rows[2,7] = 9
$rows[3,0] = 3;  $rows[3,1] = "014208"; $rows[3,2] = "华安产业精选混合C";         $rows[3,3] = "7.93";  $rows[3,4] = "62.03"; $rows[3,5] = "2.03"; $rows[3,6] = "0.1610"; $rows[3,7] = 7
$rows[4,0] = 4;  $rows[4,1] = "519656"; $rows[4,2] = "银河灵活配置混合 - A";      $rows[4,3] = "0.72";  $rows[4,4] = "59.27"; $rows[4,5] = "3.57"; $rows[4,6] = "0.0257"; $rows[4,7] = 6
$rows[5,0] = 5;  $rows[5,1] = "006181"; $rows[5,2] = "格林伯锐灵活配置混合A";     $rows[5,3] = "0.29";  $rows[5,4] = "89.68"; $rows[5,5] = "6.20"; $rows[5,6] = "0.0180"; $rows[5,7] = 2
$rows[6,0] = 6;  $rows[6,1] = "519657"; $rows[6,2] = "银河灵活配置混合 - C";      $rows[6,3] = "0.33";  $rows[6,4] = "59.27"; $rows[6,5] = "3.57"; $rows[6,6] = "0.0118"; $rows[6,7] = 6
$rows[7,0] = 7;  $rows[7,1] = "006182"; $rows[7,2] = "格林伯锐灵活配置混合C";     $rows[7,3] = "0.12";  $rows[7,4] = "89.68"; $rows[7,5] = "6.20"; $rows[7,6] = "0.0074"; $rows[7,7] = 2
$rows[8,0] = 8;  $rows[8,1] = "001899"; $rows[8,2] = "东海中证社会发展安全产业主题指数"; $rows[8,3] = "0.21";  $rows[8,4] = "90.30"; $rows[8,5] = "1.95"; $rows[8,6] = "0.0041"; $rows[8,7] = 8

$newSheet.Range("A2:H10").Value = $rows

# --- Step 2: add the 2022-Q1 row to the "总计" summary sheet ---------------

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A2:D2"))

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 2.08
$totalSheet.Range("A2").Value = 0

# Renumber the index column for the rows that shifted down one position.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

# Restore the original active sheet/selection (the edit itself doesn't
# intend to change which tab the user is looking at).
$wb.Worksheets.Item("2020-Q4").Activate()
